$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 723.5
$ws.Range("J55").Value = 1447.5
$ws.Range("L55").Value = 1447.5
$ws.Range("N55").Value = -1875.5

$ws.Range("H62").Value = 4046.25
$ws.Range("I62").Value = 4028.3333
$ws.Range("K62").Value = 4028.3333
$ws.Range("M62").Value = -3404.3333

$ws.Range("H65").Value = 4046.25
$ws.Range("I65").Value = 4028.3333
$ws.Range("K65").Value = 20141.6665
$ws.Range("M65").Value = -17021.6665

$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1947.75
$ws.Range("I2").Value = 1947.75
$ws.Range("K2").Value = 1947.75
$ws.Range("M2").Value = -1834.75

$ws.Range("H32").Value = 9549.429
$ws.Range("I32").Value = 7539.4
$ws.Range("K32").Value = 7539.4
$ws.Range("M32").Value = -7252.4

$ws.Range("H116").Value = 1947.75
$ws.Range("I116").Value = 1947.75
$ws.Range("K116").Value = 1947.75
$ws.Range("M116").Value = 346.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1947.75
$ws.Range("I3").Value = 1947.75
$ws.Range("K3").Value = 1947.75
$ws.Range("M3").Value = -1833.75

$ws.Range("H99").Value = 20998.2
$ws.Range("I99").Value = 20998.2
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 20998.2
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -19500.2
$ws.Range("N99").Value = ""

$ws.Range("H132").Value = 70000
$ws.Range("J132").Value = 70000
$ws.Range("L132").Value = 70000
$ws.Range("N132").Value = -80120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 846.75
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 793.5
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 793.5
$ws.Range("M22").Value = -550
$ws.Range("N22").Value = -1493.5

$ws.Range("H86").Value = 6226.8335
$ws.Range("I86").Value = 4198.778
$ws.Range("K86").Value = 4198.778
$ws.Range("M86").Value = -3075.778

$ws.Range("H89").Value = 6226.8335
$ws.Range("I89").Value = 4198.778
$ws.Range("K89").Value = 20993.89
$ws.Range("M89").Value = -15377.89

$ws.Range("H99").Value = 3908.3333
$ws.Range("I99").Value = 3909
$ws.Range("K99").Value = 3909
$ws.Range("M99").Value = -2411

$ws.Range("H126").Value = 3908.3333
$ws.Range("I126").Value = 3909
$ws.Range("K126").Value = 11727
$ws.Range("M126").Value = -9257

$ws.Range("H141").Value = 483195
$ws.Range("J141").Value = 698658.3
$ws.Range("L141").Value = 698658.3
$ws.Range("N141").Value = -709018.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 296.33334
$ws.Range("I26").Value = 300
$ws.Range("J26").Value = 289
$ws.Range("K26").Value = 900
$ws.Range("L26").Value = 867
$ws.Range("M26").Value = -612
$ws.Range("N26").Value = -1443

$ws.Range("H93").Value = 6666
$ws.Range("I93").Value = 4999
$ws.Range("K93").Value = 14997
$ws.Range("M93").Value = -13125

$ws.Range("H94").Value = 10000
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").Value = ""

$ws.Range("H96").Value = 5000
$ws.Range("J96").Value = 5000
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -19118

$ws.Range("H97").Value = 1259.3077
$ws.Range("J97").Value = 825.3333
$ws.Range("L97").Value = 2475.9999
$ws.Range("N97").Value = -3467.9999

$ws.Range("H98").Value = 3294.8333
$ws.Range("J98").Value = 2995
$ws.Range("L98").Value = 8985
$ws.Range("N98").Value = -11981

$ws.Range("H99").Value = 4332.778
$ws.Range("I99").Value = 1261.25
$ws.Range("K99").Value = 3783.75
$ws.Range("M99").Value = -1537.75

$ws.Range("H100").Value = 10028
$ws.Range("J100").Value = 10028
$ws.Range("L100").Value = 30084
$ws.Range("N100").Value = -31706

$ws.Range("H101").Value = 7500
$ws.Range("J101").Value = 7500
$ws.Range("L101").Value = 22500
$ws.Range("N101").Value = -27368

$ws.Range("H103").Value = 5047.5
$ws.Range("I103").Value = 95
$ws.Range("J103").Value = 10000
$ws.Range("K103").Value = 285
$ws.Range("L103").Value = 30000
$ws.Range("M103").Value = 594
$ws.Range("N103").Value = -31758

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3899
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""

$ws.Range("H83").Value = 3899
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7214.5
$ws.Range("I22").Value = 6245
$ws.Range("J22").Value = 7699.25
$ws.Range("K22").Value = 6245
$ws.Range("L22").Value = 7699.25
$ws.Range("M22").Value = -5950
$ws.Range("N22").Value = -8289.25

$ws.Range("H27").Value = 7214.5
$ws.Range("I27").Value = 6245
$ws.Range("J27").Value = 7699.25
$ws.Range("K27").Value = 6245
$ws.Range("L27").Value = 7699.25
$ws.Range("M27").Value = -6138
$ws.Range("N27").Value = -7913.25

$ws.Range("H55").Value = 858.6
$ws.Range("I55").Value = 543.4
$ws.Range("J55").Value = 1489
$ws.Range("K55").Value = 543.4
$ws.Range("L55").Value = 1489
$ws.Range("M55").Value = -370.4
$ws.Range("N55").Value = -1835

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2638.3333
$ws.Range("I4").Value = 166
$ws.Range("J4").Value = 15000
$ws.Range("K4").Value = 166
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = -53
$ws.Range("N4").Value = -15226

$ws.Range("H62").Value = 6424.375
$ws.Range("I62").Value = 2650.5
$ws.Range("K62").Value = 2650.5
$ws.Range("M62").Value = -2026.5

$ws.Range("H65").Value = 6424.375
$ws.Range("I65").Value = 2650.5
$ws.Range("K65").Value = 13252.5
$ws.Range("M65").Value = -10132.5
